# Update NATMI LR-pair results (Vcan-Itgb1) with recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2364713333333333
$ws.Range("H2").Value = 0.709414
$ws.Range("I2").Value = 0.002249544876489787
$ws.Range("J2").Value = 0.002249544876489787
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 14.43459003845178
$ws.Range("R2").Value = 129.911310346066
$ws.Range("S2").Value = 0.0004597200189761008
$ws.Range("T2").Value = 0.0004597200189761008

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2364713333333333
$ws.Range("H3").Value = 0.709414
$ws.Range("I3").Value = 0.002249544876489787
$ws.Range("J3").Value = 0.002249544876489787
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 25.14032352764133
$ws.Range("R3").Value = 226.262911748772
$ws.Range("S3").Value = 0.0008006815557909823
$ws.Range("T3").Value = 0.0008006815557909824

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2364713333333333
$ws.Range("H4").Value = 0.709414
$ws.Range("I4").Value = 0.002249544876489787
$ws.Range("J4").Value = 0.002249544876489787
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 31.05776877293244
$ws.Range("R4").Value = 279.519918956392
$ws.Range("S4").Value = 0.0009891433017227044
$ws.Range("T4").Value = 0.0009891433017227044

$ws.Range("I5").Value = 0.9862688099613843
$ws.Range("J5").Value = 0.9862688099613843
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 6328.562763201635
$ws.Range("R5").Value = 56957.06486881472
$ws.Range("S5").Value = 0.2015552215781912
$ws.Range("T5").Value = 0.2015552215781912

$ws.Range("I6").Value = 0.9862688099613843
$ws.Range("J6").Value = 0.9862688099613843
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.3510431169616131
$ws.Range("T6").Value = 0.3510431169616131

$ws.Range("I7").Value = 0.9862688099613843
$ws.Range("J7").Value = 0.9862688099613843
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 13616.66929513908
$ws.Range("R7").Value = 122550.0236562517
$ws.Range("S7").Value = 0.43367047142158
$ws.Range("T7").Value = 0.4336704714215801

$ws.Range("G8").Value = 1.206946333333333
$ws.Range("H8").Value = 3.620839
$ws.Range("I8").Value = 0.01148164516212593
$ws.Range("J8").Value = 0.01148164516212593
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 73.67394294479347
$ws.Range("R8").Value = 663.0654865031412
$ws.Range("S8").Value = 0.002346404460286103
$ws.Range("T8").Value = 0.002346404460286104

$ws.Range("G9").Value = 1.206946333333333
$ws.Range("H9").Value = 3.620839
$ws.Range("I9").Value = 0.01148164516212593
$ws.Range("J9").Value = 0.01148164516212593
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 128.3158549189913
$ws.Range("R9").Value = 1154.842694270922
$ws.Range("S9").Value = 0.004086667311032295
$ws.Range("T9").Value = 0.004086667311032296

$ws.Range("G10").Value = 1.206946333333333
$ws.Range("H10").Value = 3.620839
$ws.Range("I10").Value = 0.01148164516212593
$ws.Range("J10").Value = 0.01148164516212593
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 158.5184115706991
$ws.Range("R10").Value = 1426.665704136292
$ws.Range("S10").Value = 0.005048573390807533
$ws.Range("T10").Value = 0.005048573390807534
